# Generate Report for Handback
# Adds the missing "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" data for the
# 79152f12-fae5-4d6e-85ff-b6c42e2fc2b2 row (row 7) on both locale sheets,
# and widens the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/03acd099bb4dad81299bd8cb7834d86778e2dde3/e2e/79152f12-fae5-4d6e-85ff-b6c42e2fc2b2.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e3e3d1f0c0b5c074281d7c9dc16d7f86aa87fa68/e2e/79152f12-fae5-4d6e-85ff-b6c42e2fc2b2.md."
$handbackFileName = "79152f12-fae5-4d6e-85ff-b6c42e2fc2b2.md"
$handbackFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/03acd099bb4dad81299bd8cb7834d86778e2dde3/e2e/79152f12-fae5-4d6e-85ff-b6c42e2fc2b2.md"

$sheetInfo = @{
    "zh-cn" = "2016-08-21 04:49:44"
    "de-de" = "2016-08-21 04:49:50"
}

foreach ($sheetName in $sheetInfo.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    # I7 = "Latest Target File" -> the handback markdown file, as a hyperlink
    $ws.Range("I7").Value = $handbackFileName
    $ws.Hyperlinks.Add($ws.Range("I7"), $handbackFileUrl, "", "", $handbackFileName) | Out-Null

    # J7 = "Latest Handback File" -> same xlf file referenced by "Latest Target File" (G7)
    $targetFile = $ws.Range("G7").Value2
    $ws.Range("J7").Value = $targetFile

    # K7 = "Latest Handback DateTime"
    $ws.Range("K7").Value = $sheetInfo[$sheetName]

    # P7 = "Error Detail"
    $ws.Range("P7").Value = $errorMessage

    # Widen the Error Detail column (P / column 16) to fit the message.
    # 39.15 compensates for the engine's pixel-width rounding so the
    # stored column width comes out to exactly 40.
    $ws.Columns.Item(16).ColumnWidth = 39.15
}
